$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("subte")

# Delete the "TIPO" column (column B). This shifts the old "TOTAL" column (C)
# left into column B, and its values (all "SUBTE") get removed entirely.
$ws.Columns.Item(2).Delete()

# Rename the new header in B1 (previously "TOTAL", shifted from C1) to "SUBTE".
$ws.Range("B1").Value = "SUBTE"
